$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2-4: rokopis (manuscript) name corrections ---
$ws.Range("B2").Value = "Berke "
$ws.Range("B3").Value = "Berke "
$ws.Range("B4").Value = "NUK "

# --- Row 5: Akos Doncec / Martjanska II now has folio/status/napotnica info ---
$ws.Range("C5").Value = "celotna"
$ws.Range("D5").Value = "končano"
$ws.Range("D2").Copy()
$ws.Range("D5").PasteSpecial(-4122)
$ws.Range("E5").Value = "ne"

# --- Row 7: Kosijeva pesmarica now has folio/status/napotnica info ---
$ws.Range("C7").Value = "celotna"
$ws.Range("D7").Value = "končano"
$ws.Range("D2").Copy()
$ws.Range("D7").PasteSpecial(-4122)
$ws.Range("E7").Value = "da"

# --- Row 12: prepisovalec correction + status/napotnica added ---
$ws.Range("A12").Value = "Špeka Kovačič"
$ws.Range("D12").Value = "končano"
$ws.Range("D2").Copy()
$ws.Range("D12").PasteSpecial(-4122)
$ws.Range("E12").Value = "ne"

# --- Row 13: folio/status/napotnica added ---
$ws.Range("C13").Value = "celotna"
$ws.Range("D13").Value = "končano"
$ws.Range("D2").Copy()
$ws.Range("D13").PasteSpecial(-4122)
$ws.Range("E13").Value = "ne"

# --- Row 14: new manuscript entry ---
$ws.Range("A14").Value = "Špeka Kovačič"
$ws.Range("B14").Value = "Cantiones mortualis"
$ws.Range("C14").Value = "celotna"
$ws.Range("D14").Value = "končano"
$ws.Range("D2").Copy()
$ws.Range("D14").PasteSpecial(-4122)
$ws.Range("E14").Value = "ne"

# --- Row 15: new manuscript entry (Gaberjeva I) ---
$ws.Range("A15").Value = "Špela Kovačič"
$ws.Range("B15").Value = "Gaberjeva I"
$ws.Range("C15").Value = "CELOTNA"
$ws.Range("D15").Value = "končano"
$ws.Range("D2").Copy()
$ws.Range("D15").PasteSpecial(-4122)
$ws.Range("E15").Value = "ne"

# --- Row 16: new manuscript entry (Gaber-Bokan) ---
$ws.Range("A16").Value = "Špela Kovačič"
$ws.Range("B16").Value = "Gaber-Bokan"
$ws.Range("C16").Value = "celotna"
$ws.Range("D16").Value = "končano"
$ws.Range("D2").Copy()
$ws.Range("D16").PasteSpecial(-4122)
$ws.Range("E16").Value = "ne"

# --- Rows 17-21: new prepisovalci (transcribers) joining list ---
$ws.Range("A17").Value = "Sara Gomboc"
$ws.Range("A18").Value = "Žana Horvat"
$ws.Range("A19").Value = "Tina Raj"
$ws.Range("A20").Value = "Laura Sobočan"
$ws.Range("A21").Value = "Katja Huber"

$ws.Range("A22").Select() | Out-Null
